# Modify the original data table style: rotate columns B, C, D, E
# (both header text and the numeric data below them) for rows 1-11,
# i.e. new B = old D, new C = old B, new D = old E, new E = old C.
# Columns A and F are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 1
$lastRow = 11

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Read the original values of columns B, C, D, E first, then
    # write them back in rotated order:
    #   new B = old D
    #   new C = old B
    #   new D = old E
    #   new E = old C
    $valB = $ws.Cells.Item($r, 2).Value2
    $valC = $ws.Cells.Item($r, 3).Value2
    $valD = $ws.Cells.Item($r, 4).Value2
    $valE = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 2).Value2 = $valD
    $ws.Cells.Item($r, 3).Value2 = $valB
    $ws.Cells.Item($r, 4).Value2 = $valE
    $ws.Cells.Item($r, 5).Value2 = $valC
}
